# Update "想去人数" (interest counts) in column F across the three sheets
# that carry event rows: 展览 (Exhibitions), 本地生活 (Local Life), and
# 全部类型 (All Types - the combined view). 演出 (Performances) is untouched.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$ws1.Range("F2").Value = 70
$ws1.Range("F4").Value = 1763
$ws1.Range("F6").Value = 626
$ws1.Range("F7").Value = 1145
$ws1.Range("F8").Value = 1578
$ws1.Range("F9").Value = 172
$ws1.Range("F10").Value = 172
$ws1.Range("F12").Value = 1503
$ws1.Range("F13").Value = 3126
$ws1.Range("F14").Value = 673
$ws1.Range("F15").Value = 1818
$ws1.Range("F16").Value = 1821
$ws1.Range("F17").Value = 881
$ws1.Range("F18").Value = 297
$ws1.Range("F20").Value = 1506
$ws1.Range("F21").Value = 303
$ws1.Range("F24").Value = 1281
$ws1.Range("F26").Value = 488
$ws1.Range("F27").Value = 179
$ws1.Range("F28").Value = 5858
$ws1.Range("F29").Value = 5355
$ws1.Range("F31").Value = 592
$ws1.Range("F32").Value = 1701
$ws1.Range("F33").Value = 91
$ws1.Range("F34").Value = 221

# --- 本地生活 (sheet3) ---
$ws3.Range("F2").Value = 49

# --- 全部类型 (sheet4) ---
$ws4.Range("F2").Value = 70
$ws4.Range("F3").Value = 49
$ws4.Range("F7").Value = 1763
$ws4.Range("F9").Value = 626
$ws4.Range("F10").Value = 1145
$ws4.Range("F11").Value = 1578
$ws4.Range("F12").Value = 172
$ws4.Range("F13").Value = 172
$ws4.Range("F16").Value = 1503
$ws4.Range("F17").Value = 3126
$ws4.Range("F18").Value = 673
$ws4.Range("F19").Value = 1818
$ws4.Range("F21").Value = 881
$ws4.Range("F22").Value = 297
$ws4.Range("F24").Value = 1506
$ws4.Range("F25").Value = 303
$ws4.Range("F30").Value = 1281
$ws4.Range("F32").Value = 488
$ws4.Range("F33").Value = 179
$ws4.Range("F34").Value = 5858
$ws4.Range("F35").Value = 5355
$ws4.Range("F37").Value = 592
$ws4.Range("F38").Value = 1701
$ws4.Range("F41").Value = 91
$ws4.Range("F42").Value = 221
